$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 102, pushing existing rows 102:150 down to 103:151
$ws.Rows.Item(102).Insert()

# Populate the newly inserted row 102 with the new record
$ws.Cells.Item(102, 1).Value = 10
$ws.Cells.Item(102, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(102, 3).Value = "La Araucanía"
$ws.Cells.Item(102, 4).Value = 45029
$ws.Cells.Item(102, 4).NumberFormat = $ws.Cells.Item(103, 4).NumberFormat
$ws.Cells.Item(102, 5).Value = 9
$ws.Cells.Item(102, 6).Value = 100112035
$ws.Cells.Item(102, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(102, 8).Value = "Sin especificar"
$ws.Cells.Item(102, 9).Value = "Primera"
$ws.Cells.Item(102, 10).Value = 50
$ws.Cells.Item(102, 11).Value = 30000
$ws.Cells.Item(102, 12).Value = 30000
$ws.Cells.Item(102, 13).Value = 30000
$ws.Cells.Item(102, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(102, 15).Value = "Región Metropolitana"
$ws.Cells.Item(102, 16).Value = 2000
$ws.Cells.Item(102, 17).Value = 15
$ws.Cells.Item(102, 18).Value = "Hortaliza"
